$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 18121.334
$ws.Range("I9").Value = 25432
$ws.Range("J9").Value = 3500
$ws.Range("K9").Value = 25432
$ws.Range("L9").Value = 3500
$ws.Range("M9").Value = -25263
$ws.Range("N9").Value = -3838
$ws.Range("H32").Value = 3312.25
$ws.Range("I32").Value = 750
$ws.Range("J32").Value = 4166.3335
$ws.Range("K32").Value = 750
$ws.Range("L32").Value = 4166.3335
$ws.Range("N32").Value = -4818.3335
$ws.Range("M32").Value = -424
$ws.Range("H43").Value = 529204.5
$ws.Range("I43").Value = 2832.5557
$ws.Range("K43").Value = 2832.5557
$ws.Range("M43").Value = -2763.5557
$ws.Range("H106").Value = 1350.625
$ws.Range("I106").Value = 1390.6666
$ws.Range("J106").Value = 750
$ws.Range("K106").Value = 1390.6666
$ws.Range("L106").Value = 750
$ws.Range("M106").Value = -759.6666
$ws.Range("N106").Value = -2012
$ws.Range("H130").Value = 74599.664
$ws.Range("I130").Value = 62349.5
$ws.Range("K130").Value = 62349.5
$ws.Range("M130").Value = -57329.5
$ws.Range("H137").Value = 2548.0356
$ws.Range("J137").Value = 1884.4667
$ws.Range("L137").Value = 5653.4001
$ws.Range("N137").Value = -10753.4001
$ws.Range("H138").Value = 7412076.5
$ws.Range("I138").Value = 1188.125
$ws.Range("J138").Value = 11500843
$ws.Range("K138").Value = 3564.375
$ws.Range("L138").Value = 34502529
$ws.Range("M138").Value = 1575.625
$ws.Range("N138").Value = -34512809

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8286.27
$ws.Range("I32").Value = 7539.826
$ws.Range("J32").Value = 14009
$ws.Range("K32").Value = 7539.826
$ws.Range("L32").Value = 14009
$ws.Range("M32").Value = -7252.826
$ws.Range("N32").Value = -14583
$ws.Range("H61").Value = 29417688
$ws.Range("I61").Value = 41671224
$ws.Range("K61").Value = 41671224
$ws.Range("M61").Value = -41671012
$ws.Range("H122").Value = 3072.5938
$ws.Range("I122").Value = 2203.611
$ws.Range("J122").Value = 4189.857
$ws.Range("K122").Value = 6610.833
$ws.Range("L122").Value = 12569.571
$ws.Range("M122").Value = -4160.833
$ws.Range("N122").Value = -17469.571
$ws.Range("H132").Value = 34540884
$ws.Range("I132").Value = 11708.782
$ws.Range("K132").Value = 35126.346
$ws.Range("M132").Value = -32596.346
$ws.Range("H136").Value = 29417688
$ws.Range("I136").Value = 41671224
$ws.Range("K136").Value = 125013672
$ws.Range("M136").Value = -125011122

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1913.6666
$ws.Range("I20").Value = 1947.8
$ws.Range("K20").Value = 1947.8
$ws.Range("M20").Value = -1700.8
$ws.Range("H76").Value = 24499.5
$ws.Range("J76").Value = 19000
$ws.Range("L76").Value = 19000
$ws.Range("N76").Value = -19630
$ws.Range("H79").Value = 24499.5
$ws.Range("J79").Value = 19000
$ws.Range("L79").Value = 19000
$ws.Range("N79").Value = -21184
$ws.Range("H87").Value = 85083.336
$ws.Range("I87").Value = 87750
$ws.Range("K87").Value = 87750
$ws.Range("M87").Value = -86502
$ws.Range("H90").Value = 85083.336
$ws.Range("I90").Value = 87750
$ws.Range("K90").Value = 263250
$ws.Range("M90").Value = -257010
$ws.Range("H130").Value = 80000
$ws.Range("J130").Value = 80000
$ws.Range("L130").Value = 80000
$ws.Range("N130").Value = -90040
$ws.Range("H134").Value = 2410.9355
$ws.Range("I134").Value = 1920.6086
$ws.Range("K134").Value = 5761.825800000001
$ws.Range("M134").Value = -3226.825800000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4173.8335
$ws.Range("I31").Value = 3084.5417
$ws.Range("J31").Value = 5626.222
$ws.Range("K31").Value = 3084.5417
$ws.Range("L31").Value = 5626.222
$ws.Range("M31").Value = -2789.5417
$ws.Range("N31").Value = -6216.222
$ws.Range("H34").Value = 4173.8335
$ws.Range("I34").Value = 3084.5417
$ws.Range("J34").Value = 5626.222
$ws.Range("K34").Value = 3084.5417
$ws.Range("L34").Value = 5626.222
$ws.Range("M34").Value = -2882.5417
$ws.Range("N34").Value = -6030.222
$ws.Range("H132").Value = 63994.273
$ws.Range("I132").Value = 67219.74000000001
$ws.Range("K132").Value = 201659.22
$ws.Range("M132").Value = -199129.22

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 98.5
$ws.Range("I20").Value = 98.5
$ws.Range("K20").Value = 295.5
$ws.Range("M20").Value = -68.5
$ws.Range("H21").Value = 3673.6667
$ws.Range("J21").Value = 5500
$ws.Range("L21").Value = 16500
$ws.Range("N21").Value = -16846
$ws.Range("H39").Value = 331.5
$ws.Range("I39").Value = 304.89474
$ws.Range("K39").Value = 914.6842200000001
$ws.Range("M39").Value = -620.6842200000001
$ws.Range("H130").Value = 2665.2
$ws.Range("J130").Value = 3088.6667
$ws.Range("L130").Value = 9266.000100000001
$ws.Range("N130").Value = -19306.0001
$ws.Range("H131").Value = 1677.525
$ws.Range("J131").Value = 1718.0834
$ws.Range("L131").Value = 5154.2502
$ws.Range("N131").Value = -15234.2502
$ws.Range("H134").Value = 7195.9
$ws.Range("I134").Value = 2160
$ws.Range("J134").Value = 14749.75
$ws.Range("K134").Value = 6480
$ws.Range("L134").Value = 44249.25
$ws.Range("M134").Value = -1410
$ws.Range("N134").Value = -54389.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 80000
$ws.Range("I53").Value = 80000
$ws.Range("K53").Value = 80000
$ws.Range("M53").Value = -79369
$ws.Range("H58").Value = 25000
$ws.Range("I58").Value = 25000
$ws.Range("K58").Value = 25000
$ws.Range("M58").Value = -24723
$ws.Range("H113").Value = 3477
$ws.Range("I113").Value = 2643.8333
$ws.Range("J113").Value = 3861.5386
$ws.Range("K113").Value = 2643.8333
$ws.Range("L113").Value = 3861.5386
$ws.Range("M113").Value = -473.8332999999998
$ws.Range("N113").Value = -8201.5386
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 89333.336
$ws.Range("J131").Value = 89333.336
$ws.Range("L131").Value = 89333.336
$ws.Range("N131").Value = -99413.336
$ws.Range("H132").Value = 1479.2222
$ws.Range("I132").Value = 1479.2222
$ws.Range("K132").Value = 4437.6666
$ws.Range("M132").Value = -1907.6666
$ws.Range("H136").Value = 19818
$ws.Range("J136").Value = 19818
$ws.Range("L136").Value = 59454
$ws.Range("N136").Value = -64554

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 509.04544
$ws.Range("I55").Value = 184.33333
$ws.Range("K55").Value = 184.33333
$ws.Range("M55").Value = -11.33332999999999
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H63").Value = 67666
$ws.Range("J63").Value = 63998
$ws.Range("L63").Value = 63998
$ws.Range("N63").Value = -65496
$ws.Range("H66").Value = 67666
$ws.Range("J66").Value = 63998
$ws.Range("L66").Value = 191994
$ws.Range("N66").Value = -199482
$ws.Range("H68").Value = 6138.1665
$ws.Range("I68").Value = 2789.5
$ws.Range("J68").Value = 7812.5
$ws.Range("K68").Value = 2789.5
$ws.Range("L68").Value = 7812.5
$ws.Range("M68").Value = -2040.5
$ws.Range("N68").Value = -9310.5
$ws.Range("H71").Value = 6138.1665
$ws.Range("I71").Value = 2789.5
$ws.Range("J71").Value = 7812.5
$ws.Range("K71").Value = 13947.5
$ws.Range("L71").Value = 39062.5
$ws.Range("M71").Value = -10203.5
$ws.Range("N71").Value = -46550.5
$ws.Range("H132").Value = 7101.2354
$ws.Range("I132").Value = 4969.6665
$ws.Range("K132").Value = 14908.9995
$ws.Range("M132").Value = -12378.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 12003
$ws.Range("J18").Value = 12003
$ws.Range("L18").Value = 12003
$ws.Range("N18").Value = -12349
